# Applies the cryptos.xlsx data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.184.51'
$ws.Range("E2").Value = '  +3.40%  '

$ws.Range("D3").Value = '2.312.04'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.536'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.95%  '

$ws.Range("D9").Value = '2.329.92'
$ws.Range("E9").Value = '  +2.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.23%  '

$ws.Range("E11").Value = '  +0.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.85%  '

$ws.Range("E13").Value = '  +1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.34%  '

$ws.Range("D15").Value = '2.725.63'
$ws.Range("E15").Value = '  +2.08%  '

$ws.Range("D16").Value = '56.355.60'
$ws.Range("E16").Value = '  +3.77%  '

$ws.Range("E17").Value = '  +4.49%  '

$ws.Range("D18").Value = '2.310.50'
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.87%  '

$ws.Range("E20").Value = '  +3.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.158'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.11%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.65%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0731'
$ws.Range("E30").Value = '  +5.70%  '

$ws.Range("E31").Value = '  +10.06%  '

$ws.Range("E32").Value = '  +4.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.82%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.991'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.55%  '

$ws.Range("E36").Value = '  +5.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.924'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.09%  '

$ws.Range("E39").Value = '  +8.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.50'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.49%  '

$ws.Range("E41").Value = '  +2.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '276.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0510'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0927'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.557'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.382'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0216'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.70%  '
